{"js": "// 1) Title on the first page: \"ICLR 2022\" -> \"ICLR 2023\"\nconst title = context.document.body.search(\"ICLR 2022\", { matchCase: true });\ntitle.load(\"items,text\");\nawait context.sync();\nif (title.items.length > 0) {\n  title.items[0].insertText(\"ICLR 2023\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Inline mention \"ICLR2022\" -> \"ICLR2023\" (immediately followed by a\n//    _GoBack bookmark, as Word stamps after the most recent edit point).\nconst inline = context.document.body.search(\"ICLR2022\", { matchCase: true });\ninline.load(\"items,text\");\nawait context.sync();\nif (inline.items.length > 0) {\n  const run = inline.items[0];\n  run.insertText(\"ICLR2023\", \"Replace\");\n  await context.sync();\n\n  const after = run.getRange(\"After\");\n  after.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 3) The workshop-website paragraph: drop the hyperlink entirely and\n//    replace its visible text with the new URL (plain text, no link).\nconst link = context.document.body.search(\"https://africanlp.masakhane.io/\", { matchCase: true });\nlink.load(\"items,text\");\nawait context.sync();\nif (link.items.length > 0) {\n  const linkRange = link.items[0];\n\n  // Remove the hyperlink field, leaving a normal run behind.\n  linkRange.hyperlink = \"\";\n  await context.sync();\n\n  // Re-search (the range above may be stale after the hyperlink edit).\n  const link2 = context.document.body.search(\"https://africanlp.masakhane.io/\", { matchCase: true });\n  link2.load(\"items,text\");\n  await context.sync();\n  const run2 = link2.items[0];\n  run2.insertText(\"https://sites.google.com/view/africanlp2023\", \"Replace\");\n  await context.sync();\n\n  // Normalize the paragraph indentation to match the retyped line.\n  const para = run2.paragraphs.getFirst();\n  para.leftIndent = 0;\n  para.firstLineIndent = 0;\n  await context.sync();\n}\n\n// 4) Footer \"Under review: AfricaNLP workshop at ICLR2022\" -> ICLR2023.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const header = section.getHeader(\"Primary\");\n  const headerHits = header.search(\"ICLR2022\", { matchCase: true });\n  headerHits.load(\"items,text\");\n  await context.sync();\n  for (const hit of headerHits.items) {\n    hit.insertText(\"ICLR2023\", \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title on the first page: \"ICLR 2022\" -> \"ICLR 2023\"\n$titleRange = $d.Content\nif ($titleRange.Find.Execute(\"ICLR 2022\")) {\n    $titleRange.Text = \"ICLR 2023\"\n}\n\n# 2) Inline mention \"ICLR2022\" -> \"ICLR2023\", followed immediately by an\n#    empty \"_GoBack\" bookmark (Word always drops one at the last edit\n#    point). We briefly insert a one-character placeholder after the new\n#    text so the bookmark position isn't the very last character before a\n#    paragraph mark, add the bookmark there, then remove the placeholder -\n#    the bookmark collapses back to a true zero-length mark.\n$inlineRange = $d.Content\nif ($inlineRange.Find.Execute(\"ICLR2022\")) {\n    $inlineRange.Text = \"ICLR2023\"\n\n    $afterRange = $inlineRange.Duplicate\n    $afterRange.Collapse(0)\n    $afterRange.InsertAfter(\"X\")\n\n    $bmPos = $afterRange.Start\n    $bmRange = $d.Range($bmPos, $bmPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n    $placeholder = $d.Range($bmPos, $bmPos + 1)\n    $placeholder.Delete()\n}\n\n# 3) Workshop-website paragraph: drop the hyperlink field (keep the run)\n#    and retype its text with the new URL, then square up the paragraph\n#    indentation to match the retyped line.\n$hyperlinks = $d.Hyperlinks\nfor ($i = 1; $i -le $hyperlinks.Count; $i++) {\n    $h = $hyperlinks.Item($i)\n    if ($h.Address -eq \"https://africanlp.masakhane.io/\") {\n        $h.Address = \"\"\n        break\n    }\n}\n\n$linkTextRange = $d.Content\nif ($linkTextRange.Find.Execute(\"https://africanlp.masakhane.io/\")) {\n    $linkTextRange.Text = \"https://sites.google.com/view/africanlp2023\"\n    $linkPara = $linkTextRange.Paragraphs.Item(1)\n    $linkPara.Range.ParagraphFormat.LeftIndent = 0\n    $linkPara.Range.ParagraphFormat.FirstLineIndent = 0\n}\n\n# 4) Default (Primary) page header: \"...workshop at ICLR2022\" -> ICLR2023.\n$section = $d.Sections.Item(1)\n$primaryHeader = $section.Headers.Item(1)\n$headerRange = $primaryHeader.Range\nif ($headerRange.Find.Execute(\"ICLR2022\")) {\n    $headerRange.Text = \"ICLR2023\"\n}\n"}
